# Update the "avatar" column (C) for every employee row so that it points
# at the same Baidu-hosted image URL, and drop the stray hyperlink that
# used to live on C3 (it pointed at a different image than its display
# text). Finally move the active selection to C13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baiduUrl = "https://img1.baidu.com/it/u=2165937980,813753762&fm=253&fmt=auto&app=138&f=JPEG?w=500&h=500"

# C2 just needs its text swapped for the new URL.
$ws.Range("C2").Value = $baiduUrl

# C3 carries an actual hyperlink (different target than its displayed
# text) plus the associated hyperlink character formatting - strip both
# before writing the new, plain value.
$ws.Range("C3").Hyperlinks.Delete()
$ws.Range("C3").ClearFormats()
$ws.Range("C3").Value = $baiduUrl

# The remaining avatar cells (C4:C37) simply get the same URL.
for ($r = 4; $r -le 37; $r++) {
    $ws.Cells.Item($r, 3).Value = $baiduUrl
}

# Match the author's final cursor position.
[void]$ws.Activate()
[void]$ws.Range("C13").Select()
